$d = $word.ActiveDocument

# ============================================================
# Article 1 (The Servant as Leader) - fill in empty answer paragraphs
# ============================================================

# Q: "Have you ever worked with a leader..." -> paragraph 11
$d.Paragraphs.Item(11).Range.Text = "One of my past managers was kind of similar. He was very easy going and preferred to be seen as more of a team lead role than a manager."

# Q: "Is the present-day world too jaded..." -> paragraph 13
$d.Paragraphs.Item(13).Range.Text = "No. If this style really did work then the organizations that followed it would have an advantage over the ones that didn’t. It’s not a common leadership style because it’s not very effective. "

# Q: "Do you see in yourself any of the qualities..." -> paragraph 15
$d.Paragraphs.Item(15).Range.Text = "Primarily the part where he mentions that people should look to solve a problem instead of blaming it on others. I like that because it encourages personal responsibility instead of playing blame games."

# Q: "How can you develop or enhance these qualities?" -> paragraph 17
$d.Paragraphs.Item(17).Range.Text = "Keep them in mind when interacting with other people. Make an effort to correct yourself when you take some action that would be against the qualities mentioned."

# Q: "Are these qualities or can they be learned?" -> insert " innate" after "qualities"
$p18 = $d.Paragraphs.Item(18)
$insertPt = $d.Range($p18.Range.Start + 19, $p18.Range.Start + 19)
$insertPt.InsertAfter(" innate")

# New paragraph 19 (previously empty) gets an answer
$d.Paragraphs.Item(19).Range.Text = "As with all nature versus nurture questions, people are predisposed towards their personalities but are still capable of changing. "

Write-Host "Stage 1 complete"

# ============================================================
# Table 2 (Chaordic Leadership article) - merge "Chaordic" runs
# ============================================================
$t2 = $d.Tables.Item(2)
$cell = $t2.Cell(2, 3)
$cellPara = $cell.Range.Paragraphs.Item(1)
$cellFull = $cellPara.Range
$chaordicRng = $d.Range($cellFull.Start + 14, $cellFull.End)
$chaordicRng.Text = " Chaordic Leadership"

Write-Host "Stage 2 complete"
